# Regenerate merged AHB files
# - rename header labels "_old" -> "_FV2304" and "_new" -> "_FV2310"
# - turn the data range into an Excel Table (ListObject) with autofilter
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row labels (row 1, columns A:J and L:U) ---
$oldToNew = @{
    "Segmentname_old"           = "Segmentname_FV2304"
    "Segmentgruppe_old"         = "Segmentgruppe_FV2304"
    "Segment_old"               = "Segment_FV2304"
    "Datenelement_old"          = "Datenelement_FV2304"
    "Segment ID_old"            = "Segment ID_FV2304"
    "Code_old"                  = "Code_FV2304"
    "Qualifier_old"             = "Qualifier_FV2304"
    "Beschreibung_old"          = "Beschreibung_FV2304"
    "Bedingungsausdruck_old"    = "Bedingungsausdruck_FV2304"
    "Bedingung_old"             = "Bedingung_FV2304"
    "Segmentname_new"           = "Segmentname_FV2310"
    "Segmentgruppe_new"         = "Segmentgruppe_FV2310"
    "Segment_new"               = "Segment_FV2310"
    "Datenelement_new"          = "Datenelement_FV2310"
    "Segment ID_new"            = "Segment ID_FV2310"
    "Code_new"                  = "Code_FV2310"
    "Qualifier_new"             = "Qualifier_FV2310"
    "Beschreibung_new"          = "Beschreibung_FV2310"
    "Bedingungsausdruck_new"    = "Bedingungsausdruck_FV2310"
    "Bedingung_new"             = "Bedingung_FV2310"
}

$headerRange = $ws.Range("A1:U1")
for ($c = 1; $c -le 21; $c++) {
    $cell = $headerRange.Cells.Item(1, $c)
    $current = [string]$cell.Value()
    if ($oldToNew.ContainsKey($current)) {
        $cell.Value = $oldToNew[$current]
    }
}

# --- 2. Convert A1:U58 into a real table (ListObject) with autofilter ---
$tableRange = $ws.Range("A1:U58")
$lo = $ws.ListObjects.Add(1, $tableRange, $false, 1)
$lo.Name = "Table1"

# --- 3. Freeze the header row (pane split after row 1) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
